$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Levantamento de horas")

# Add the hour-computation formulas to column F (rows 6-17), referencing
# the rate lookup table in columns K/L/M.
$ws.Range("F6").Formula  = "=(C6*`$K`$5)+(D6*`$L`$5)+(E6*`$M`$5)"
$ws.Range("F7").Formula  = "=(C7*`$K`$6)+(D7*`$L`$6)+(E7*`$M`$6)"
$ws.Range("F8").Formula  = "=(C8*`$K`$6)+(D8*`$L`$6)+(E8*`$M`$6)"
$ws.Range("F9").Formula  = "=(C9*`$K`$8)+(D9*`$L`$8)+(E9*`$M`$8)"
$ws.Range("F10").Formula = "=(C10*`$K`$9)+(D10*`$L`$9)+(E10*`$M`$9)"
$ws.Range("F11").Formula = "=(C11*`$K`$10)+(D11*`$L`$10)+(E11*`$M`$10)"
$ws.Range("F12").Formula = "=(C12*`$K`$11)+(D12*`$L`$11)+(E12*`$M`$11)"
$ws.Range("F13").Formula = "=(C13*`$K`$12)+(D13*`$L`$12)+(E13*`$M`$12)"
$ws.Range("F14").Formula = "=(C14*`$K`$13)+(D14*`$L`$13)+(E14*`$M`$13)"
$ws.Range("F15").Formula = "=(C15*`$K`$14)+(D15*`$L`$14)+(E15*`$M`$14)"
$ws.Range("F16").Formula = "=(C16*`$K`$15)+(D16*`$L`$15)+(E16*`$M`$15)"
$ws.Range("F17").Formula = "=(C17*`$K`$16)+(D17*`$L`$16)+(E17*`$M`$16)"

# Re-merge A2:B4 as a single merged block (was previously three separate
# merges A2:B2, A3:B3, A4:B4).
$ws.Range("A2:B2").UnMerge()
$ws.Range("A3:B3").UnMerge()
$ws.Range("A4:B4").UnMerge()
$ws.Range("A2:B4").Merge()

# Clear the stray formatted cell outside of the real data area so the
# sheet's used range (dimension) shrinks back down to A1:M19.
$ws.Range("O23").Clear()

# Update the selection on the sheet.
$ws.Range("A1:F1").Select()

$wb.Save()
